# Actualización automática hashcode mar feb 11 01:47:36 CET 2020
# Updates specific hashcode values (column B) in the hashcode.csv sheet
# while leaving the rest of the workbook untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hashcode.csv")

$ws.Range("B27").Value = "91391565d507442005cab71becb1f189"
$ws.Range("B193").Value = "96b7225b0510150b80e8851a933770f2"
$ws.Range("B273").Value = "8c6e2b75376b8490b816902250befb49"
$ws.Range("B280").Value = "492e2c8558dd97a9d0b23fd3851b3587"
$ws.Range("B297").Value = "94be80de020cde18c28d763027eccc29"
$ws.Range("B300").Value = "ceb5db6d9ed2eeaff9635a6d64dfc3bb"
$ws.Range("B358").Value = "fa7d097d702a3fb4c02ab0e6ec24568c"
$ws.Range("B367").Value = "5aae822d0f2ec57225edd7b281024f8a"
$ws.Range("B397").Value = "6669c1951a57a64b8a545ff193a79694"
$ws.Range("B416").Value = "5bc66926ec0893680b606c0d50c3c2cd"
$ws.Range("B422").Value = "9a547834bfb45447f35c0898ffcfce14"
$ws.Range("B477").Value = "3f6233748c9d480d537076a8e25cd463"
$ws.Range("B510").Value = "226ef832ff5c84d2f7ef1295940c9ce5"
$ws.Range("B511").Value = "2c7c22ed1ce7767e03ff2638310fc76b"
$ws.Range("B520").Value = "922820cb546d4143611e0ac0c6cb3e5c"
$ws.Range("B529").Value = "4725d2dc189712fda585ce4142710523"
$ws.Range("B546").Value = "74647f72a3eb673cbf036ed249f3a9f6"
$ws.Range("B564").Value = "391d814176dc4f5df88f9e713496af7e"
$ws.Range("B577").Value = "0e540c3ca72bac2a30b6f2007ef86776"
$ws.Range("B589").Value = "0311ad095aa00adb2660907f628b57db"
$ws.Range("B770").Value = "c9cf90bca9b52c9fceeea36cf5d8debc"
$ws.Range("B789").Value = "3530be274c9da14179c1054bb965cea0"
$ws.Range("B803").Value = "09b2547196d057257fa8d355bc56555a"
$ws.Range("B897").Value = "d0b7821b6e6a30385eb91b05fb4adc9f"
$ws.Range("B905").Value = "dc5ab44aaf01eeca4909629fce968836"
$ws.Range("B963").Value = "3f574683856d8cc29639b08f7ab41e07"
$ws.Range("B967").Value = "ec8951b0c90004edf34c721157014b9d"
